$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the "Results" column (column B) entirely - this shifts the
# remaining columns (Date, SIJ Tip, Standoff, Wavegen, Samples) one to the left.
$ws.Columns("B").Delete()

# The worksheet's hidden AutoFilter defined name still points at the old
# $G$2 corner; shrink it to match the new (one-column-narrower) data range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$2"
    }
}

# Leave the selection where Excel ended up after the edit.
$ws.Range("D8").Select()
